# Fruta / hortaliza, semanal
# Insert a new week of data (2 rows) at the top of the date-ordered block
# (previously rows 583-584), shifting the existing rows 583-625 down to 585-627.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 583, pushing everything from the
# old row 583 onward down by two rows.
$ws.Rows.Item(583).Insert()
$ws.Rows.Item(583).Insert()

# Populate the first new row (583) with the new weekly data.
$ws.Range("A583").Value = 8
$ws.Range("B583").Value = "Terminal La Palmera de La Serena"
$ws.Range("C583").Value = "Coquimbo"
$ws.Range("D583").Value = 45265
$ws.Range("E583").Value = 4
$ws.Range("F583").Value = 100112021
$ws.Range("G583").Value = "Ají"
$ws.Range("H583").Value = "Inferno"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 400
$ws.Range("K583").Value = 33000
$ws.Range("L583").Value = 34000
$ws.Range("M583").Value = 33500
$ws.Range("N583").Value = "$/caja 15 kilos"
$ws.Range("O583").Value = "Provincia de Limarí"
$ws.Range("P583").Value = 2233
$ws.Range("Q583").Value = 15
$ws.Range("R583").Value = "Hortaliza"

# Populate the second new row (584) with the new weekly data.
$ws.Range("A584").Value = 8
$ws.Range("B584").Value = "Terminal La Palmera de La Serena"
$ws.Range("C584").Value = "Coquimbo"
$ws.Range("D584").Value = 45265
$ws.Range("E584").Value = 4
$ws.Range("F584").Value = 100112021
$ws.Range("G584").Value = "Ají"
$ws.Range("H584").Value = "Inferno"
$ws.Range("I584").Value = "Segunda"
$ws.Range("J584").Value = 260
$ws.Range("K584").Value = 19000
$ws.Range("L584").Value = 20000
$ws.Range("M584").Value = 19500
$ws.Range("N584").Value = "$/caja 15 kilos"
$ws.Range("O584").Value = "Provincia de Limarí"
$ws.Range("P584").Value = 1300
$ws.Range("Q584").Value = 15
$ws.Range("R584").Value = "Hortaliza"
